# Daily attendance processing - swap the order of the "Recorded By" names
# in column G from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
